$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Colaboradores": add a small "Fórmulas" note block below the
# existing table (rows 15-16).
# ---------------------------------------------------------------------
$wsColab = $wb.Worksheets.Item("Colaboradores")
$wsColab.Range("A15").Value = "Fórmulas"
$wsColab.Range("A16").Value = "Salário Mensal = Salário por hora * Carga Horária mensal"

# ---------------------------------------------------------------------
# Sheet "Folha de Pagamento": add a "Formula" notes block below the
# existing table (rows 18-24).
# ---------------------------------------------------------------------
$wsFolha = $wb.Worksheets.Item("Folha de Pagamento")
$wsFolha.Range("A18").Value = "Formula"
$wsFolha.Range("A19").Value = "Horas normais = Horas Mensais * valor por hora"
$wsFolha.Range("A20").Value = "Horas Extras -> buscar total na aba de horas extras"
$wsFolha.Range("A21").Value = "INSS"
$wsFolha.Range("B21").Value = "fazer junto com professor"
$wsFolha.Range("A22").Value = "IRPF"
$wsFolha.Range("B22").Value = "fazer junto com professor"
$wsFolha.Range("A23").Value = "Outros descontos -> somar todos os descontos da aba Descontos"
$wsFolha.Range("A24").Value = "Valor a recerber = Horas normais + Horas Extras - INSS - IRPF - Outros Descontos"

# ---------------------------------------------------------------------
# Sheet "HorasExtras": insert a new "Valor por hora" column (before the
# old "Normal R$" column), format it like its neighbours, center the
# "Horas" column data, and add a "Fórmulas" notes block below the table.
# ---------------------------------------------------------------------
$wsHoras = $wb.Worksheets.Item("HorasExtras")

$wsHoras.Columns("E:E").Insert()
$wsHoras.Range("F2:F10").Copy()
$wsHoras.Range("E2:E10").PasteSpecial(-4122)
$wsHoras.Range("E2").Value = "Valor por hora"
$wsHoras.Columns("E:E").ColumnWidth = 13.7109375

$wsHoras.Range("D3:D10").HorizontalAlignment = -4108

$wsHoras.Range("A13").Value = "Fórmulas"
$wsHoras.Range("A14").Value = "Valor por hora -> buscar na planilha de colaboradores"
$wsHoras.Range("A15").Value = "Normal = Horas * valor por hora"
$wsHoras.Range("A16").Value = "Acréscimo = Normal * acrescimos sobre horas extras"
$wsHoras.Range("A17").Value = "Total = Normal + acréscimo"

# ---------------------------------------------------------------------
# Restore / update the selections on each sheet, then make "Folha de
# Pagamento" the active tab (matches activeTab="1" in workbook.xml).
# ---------------------------------------------------------------------
$wsColab.Range("A17").Select()
$wsHoras.Range("B20").Select()

$wsFolha.Activate()
$wsFolha.Range("A25").Select()
